$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.103.53"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "1.656.48"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "215.27"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "1.891.68"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("D13").Value = "1.659.82"
$ws.Range("E13").Value = "  +4.11%  "
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "65.10"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "27.092.43"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "238.44"
$ws.Range("E18").Value = "  +3.74%  "
$ws.Range("D19").Value = "7.87"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "4.42"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("D23").Value = "2.24"
$ws.Range("D24").Value = "9.23"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("D25").Value = "145.76"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("D33").Value = "1.514.44"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +10.41%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "0.576"
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").Value = "0.888"
$ws.Range("E38").Value = "  +8.53%  "
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +4.58%  "
$ws.Range("D43").Value = "66.07"
$ws.Range("E43").Value = "  +9.40%  "
$ws.Range("D44").Value = "1.798.58"
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "89.60"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "0.0977"
$ws.Range("E51").Value = "  +3.31%  "
